# Update the '想去人数' (column F) values for the rows that changed
# in the data refresh, across all four sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 376
$ws.Range("F6").Value = 771
$ws.Range("F7").Value = 212
$ws.Range("F8").Value = 1086
$ws.Range("F9").Value = 281
$ws.Range("F11").Value = 355
$ws.Range("F12").Value = 629
$ws.Range("F18").Value = 836
$ws.Range("F26").Value = 155
$ws.Range("F31").Value = 242
$ws.Range("F32").Value = 1037
$ws.Range("F33").Value = 73
$ws.Range("F35").Value = 275

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F14").Value = 590
$ws.Range("F24").Value = 294
$ws.Range("F25").Value = 263
$ws.Range("F26").Value = 3772
$ws.Range("F33").Value = 130

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 991
$ws.Range("F9").Value = 1250

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 991
$ws.Range("F7").Value = 1250
$ws.Range("F11").Value = 376
$ws.Range("F12").Value = 771
$ws.Range("F13").Value = 212
$ws.Range("F15").Value = 1086
$ws.Range("F16").Value = 281
$ws.Range("F17").Value = 355
$ws.Range("F18").Value = 629
$ws.Range("F24").Value = 836
$ws.Range("F31").Value = 155
$ws.Range("F34").Value = 590
$ws.Range("F35").Value = 590
$ws.Range("F39").Value = 242
$ws.Range("F44").Value = 294
$ws.Range("F45").Value = 294
$ws.Range("F46").Value = 263
$ws.Range("F47").Value = 1037
$ws.Range("F49").Value = 130
$ws.Range("F51").Value = 275

